$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date column C was bumped from serial date
# 45180 (2023-09-11) to 45181 (2023-09-12) for every data row (2-12).
$ws.Range("C2:C12").Value = 45181
